$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing current rows 6-8 down to 7-9
$ws.Rows.Item(6).Insert()

# Fill the new row 6 with data (same "constant" columns as surrounding rows,
# plus the specific new values from the diff)
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "Vega Modelo de Temuco"
$ws.Range("C6").Value = "La Araucanía"
$ws.Range("D6").Value = 45014
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 100112041
$ws.Range("G6").Value = "Fruto del paraíso"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 8000
$ws.Range("L6").Value = 8000
$ws.Range("M6").Value = 8000
$ws.Range("N6").Value = "`$/caja 18 kilos empedrada"
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 444
$ws.Range("Q6").Value = 18
$ws.Range("R6").Value = "Hortaliza"
